$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 39
$ws.Range("I9").Value = 26.2
$ws.Range("K9").Value = 26.2
$ws.Range("M9").Value = 142.8
$ws.Range("H53").Value = 753.625
$ws.Range("I53").Value = 676.8
$ws.Range("J53").Value = 881.6667
$ws.Range("K53").Value = 676.8
$ws.Range("L53").Value = 881.6667
$ws.Range("M53").Value = -39.79999999999995
$ws.Range("N53").Value = -2155.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 995.6667
$ws.Range("I2").Value = 998
$ws.Range("J2").Value = 991
$ws.Range("K2").Value = 998
$ws.Range("L2").Value = 991
$ws.Range("M2").Value = -885
$ws.Range("N2").Value = -1217
$ws.Range("H63").Value = 2952.6
$ws.Range("I63").Value = 2952.6
$ws.Range("K63").Value = 2952.6
$ws.Range("M63").Value = -2266.6
$ws.Range("H66").Value = 2952.6
$ws.Range("I66").Value = 2952.6
$ws.Range("K66").Value = 14763
$ws.Range("M66").Value = -11331
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H116").Value = 995.6667
$ws.Range("I116").Value = 998
$ws.Range("J116").Value = 991
$ws.Range("K116").Value = 998
$ws.Range("L116").Value = 991
$ws.Range("M116").Value = 1296
$ws.Range("N116").Value = -5579

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 995.6667
$ws.Range("I3").Value = 998
$ws.Range("J3").Value = 991
$ws.Range("K3").Value = 998
$ws.Range("L3").Value = 991
$ws.Range("M3").Value = -884
$ws.Range("N3").Value = -1219
$ws.Range("H80").Value = 1991.5
$ws.Range("I80").Value = 1994
$ws.Range("J80").Value = 1984
$ws.Range("K80").Value = 1994
$ws.Range("L80").Value = 1984
$ws.Range("M80").Value = -996
$ws.Range("N80").Value = -3980
$ws.Range("H83").Value = 1991.5
$ws.Range("I83").Value = 1994
$ws.Range("J83").Value = 1984
$ws.Range("K83").Value = 9970
$ws.Range("L83").Value = 9920
$ws.Range("M83").Value = -4978
$ws.Range("N83").Value = -19904

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1486.875
$ws.Range("I16").Value = 1297.2
$ws.Range("J16").Value = 1803
$ws.Range("K16").Value = 1297.2
$ws.Range("L16").Value = 1803
$ws.Range("M16").Value = -1010.2
$ws.Range("N16").Value = -2377
$ws.Range("H51").Value = 13994.5
$ws.Range("I51").Value = 13994.5
$ws.Range("K51").Value = 13994.5
$ws.Range("M51").Value = -13258.5
$ws.Range("H58").Value = 5499.467
$ws.Range("I58").Value = 2149.3
$ws.Range("K58").Value = 2149.3
$ws.Range("M58").Value = -1946.3
$ws.Range("H61").Value = 13994.5
$ws.Range("I61").Value = 13994.5
$ws.Range("K61").Value = 13994.5
$ws.Range("M61").Value = -13646.5
$ws.Range("H113").Value = 1486.875
$ws.Range("I113").Value = 1297.2
$ws.Range("J113").Value = 1803
$ws.Range("K113").Value = 1297.2
$ws.Range("L113").Value = 1803
$ws.Range("M113").Value = 872.8
$ws.Range("N113").Value = -6143
$ws.Range("H132").Value = 7485.5713
$ws.Range("I132").Value = 7349.75
$ws.Range("J132").Value = 7666.6665
$ws.Range("K132").Value = 22049.25
$ws.Range("L132").Value = 22999.9995
$ws.Range("M132").Value = -19519.25
$ws.Range("N132").Value = -28059.9995
$ws.Range("H136").Value = 5499.467
$ws.Range("I136").Value = 2149.3
$ws.Range("K136").Value = 6447.900000000001
$ws.Range("M136").Value = -3897.900000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 13.230769
$ws.Range("I7").Value = 10
$ws.Range("K7").Value = 30
$ws.Range("M7").Value = 82
$ws.Range("H23").Value = 185.2
$ws.Range("I23").Value = 219.6
$ws.Range("J23").Value = 150.8
$ws.Range("K23").Value = 658.8
$ws.Range("L23").Value = 452.4
$ws.Range("M23").Value = -423.8
$ws.Range("N23").Value = -922.4000000000001
$ws.Range("H80").Value = 1999.5
$ws.Range("I80").Value = 1999.5
$ws.Range("K80").Value = 5998.5
$ws.Range("M80").Value = -5062.5
$ws.Range("H83").Value = 1999.5
$ws.Range("I83").Value = 1999.5
$ws.Range("K83").Value = 17995.5
$ws.Range("M83").Value = -13315.5
$ws.Range("H92").Value = 1332.6666
$ws.Range("J92").Value = 1499
$ws.Range("L92").Value = 4497
$ws.Range("N92").Value = -6993
$ws.Range("H102").Value = 1999
$ws.Range("I102").Value = 1999
$ws.Range("K102").Value = 5997
$ws.Range("M102").Value = -3563
$ws.Range("H109").Value = 2815.125
$ws.Range("I109").Value = 3003
$ws.Range("K109").Value = 9009
$ws.Range("M109").Value = -7969

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 602.5
$ws.Range("I97").Value = 535
$ws.Range("K97").Value = 535
$ws.Range("M97").Value = -39

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1533.3334
$ws.Range("I68").Value = 1533.3334
$ws.Range("K68").Value = 1533.3334
$ws.Range("M68").Value = -784.3334
$ws.Range("H71").Value = 1533.3334
$ws.Range("I71").Value = 1533.3334
$ws.Range("K71").Value = 7666.666999999999
$ws.Range("M71").Value = -3922.666999999999
$ws.Range("H132").Value = 6991.1816
$ws.Range("I132").Value = 6426.6313
$ws.Range("J132").Value = 10566.667
$ws.Range("K132").Value = 19279.8939
$ws.Range("L132").Value = 31700.001
$ws.Range("M132").Value = -16749.8939
$ws.Range("N132").Value = -36760.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 42840.6
$ws.Range("I62").Value = 3766.6667
$ws.Range("J62").Value = 101451.5
$ws.Range("K62").Value = 3766.6667
$ws.Range("L62").Value = 101451.5
$ws.Range("M62").Value = -3142.6667
$ws.Range("N62").Value = -102699.5
$ws.Range("H65").Value = 42840.6
$ws.Range("I65").Value = 3766.6667
$ws.Range("J65").Value = 101451.5
$ws.Range("K65").Value = 18833.3335
$ws.Range("L65").Value = 507257.5
$ws.Range("M65").Value = -15713.3335
$ws.Range("N65").Value = -513497.5
$ws.Range("H113").Value = 1871.5
$ws.Range("I113").Value = 1276.4546
$ws.Range("K113").Value = 3829.3638
$ws.Range("M113").Value = -1659.3638
